{"js": "// Partie 2, Question 4 \u2014 add the new heading + answer paragraphs right\n// before the trailing (empty) paragraph that only carries the `_GoBack`\n// bookmark, i.e. right after the paragraph that ends the Question 3 answer\n// (\"...SDL_RenderPresent().\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n// items[items.length - 1] is the trailing empty \"_GoBack\" paragraph; the new\n// content must be inserted right before it.\nconst anchorParagraph = items[items.length - 2];\n\n// Insert both new paragraphs first while still plain text, then apply bold\n// only to the heading \u2014 doing it in this order keeps the bold formatting\n// from leaking onto the answer paragraph that follows.\nconst headingParagraph = anchorParagraph.insertParagraph(\"\", \"After\");\nconst answerParagraph = headingParagraph.insertParagraph(\"\", \"After\");\n\nheadingParagraph.insertText(\n  \"4) Les images anim\u00e9es : Gestion des Tileset Expliquez ce qu\\u2019est un tileset, \u00e0 quoi cela peut servir? 10 points.\",\n  \"End\"\n);\nanswerParagraph.insertText(\n  \"Un tileset est une collection d\\u2019images rectangulaires (ou bien carr\u00e9es) affich\u00e9es \u00e0 l\\u2019\u00e9cran. \" +\n    \"Le tileset correspond \u00e0 l\\u2019enti\u00e8ret\u00e9 des images utilis\u00e9es pour cr\u00e9er le rendu.\",\n  \"End\"\n);\n\nheadingParagraph.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Partie 2, Question 4 \u2014 add the new heading + answer paragraphs right\n# before the trailing (empty) paragraph that only carries the `_GoBack`\n# bookmark, i.e. right after the paragraph that ends the Question 3 answer\n# (\"...SDL_RenderPresent().\").\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$n = $paras.Count\n$anchor = $paras.Item($n - 1)\n\n# Insert both new (still plain) paragraphs first, then format the heading \u2014\n# this keeps the bold from \"leaking\" onto the answer paragraph below it.\n$anchor.Range.InsertParagraphAfter()\n$headingPara = $d.Paragraphs.Item($n)\n$headingPara.Range.InsertParagraphAfter()\n$answerPara = $d.Paragraphs.Item($n + 1)\n\n$headingPara.Range.Text = \"4) Les images anim\u00e9es : Gestion des Tileset Expliquez ce qu\u2019est un tileset, \u00e0 quoi cela peut servir? 10 points.\"\n$answerPara.Range.Text = \"Un tileset est une collection d\u2019images rectangulaires (ou bien carr\u00e9es) affich\u00e9es \u00e0 l\u2019\u00e9cran. Le tileset correspond \u00e0 l\u2019enti\u00e8ret\u00e9 des images utilis\u00e9es pour cr\u00e9er le rendu.\"\n\n$headingPara.Range.Bold = 1\n"}
